# Modify fixed O&M costs for early retirement
# Target sheet: "INS"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# E3 header: "Attrib_Cond" -> "Year"
$ws.Range("E3").Value = "Year"

# Row 4: set E4 = 2018, and change H4:AH4 from 2 to 0
$ws.Range("E4").Value = 2018
$ws.Range("H4:AH4").Value = 0

# Row 5: set E5 = 2020, H5:AH5 = 0, AJ5/AK5 labels
$ws.Range("E5").Value = 2020
$ws.Range("H5:AH5").Value = 0
$ws.Range("AJ5").Value = "T-CAR-ICE*"
$ws.Range("AK5").Value = "*Existing"

# Row 6 (new): E6 = 2025, H6:AH6 = 1, AJ6/AK6 labels
$ws.Range("E6").Value = 2025
$ws.Range("H6:AH6").Value = 1
$ws.Range("AJ6").Value = "T-CAR-ICE*"
$ws.Range("AK6").Value = "*Existing"

# Selection / view adjustments
$ws.Range("J12").Select()
